$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = New-Object 'object[,]' 24,1
$colB[0,0] = 9.327399406679785
$colB[1,0] = 8.949646982079191
$colB[2,0] = 8.709834113255562
$colB[3,0] = 8.610262670165824
$colB[4,0] = 8.593621656073479
$colB[5,0] = 8.70849854276825
$colB[6,0] = 9.19885919484997
$colB[7,0] = 10.09293208870264
$colB[8,0] = 10.70287944492555
$colB[9,0] = 10.96921257602166
$colB[10,0] = 11.06839984081079
$colB[11,0] = 11.04711319242455
$colB[12,0] = 10.97740653150596
$colB[13,0] = 10.93449026043295
$colB[14,0] = 10.68524400597144
$colB[15,0] = 10.5294373372114
$colB[16,0] = 10.43877666445965
$colB[17,0] = 10.40790325143945
$colB[18,0] = 10.54613186590998
$colB[19,0] = 10.99792679653842
$colB[20,0] = 11.2834585033547
$colB[21,0] = 11.13197566453775
$colB[22,0] = 10.53858764496759
$colB[23,0] = 9.858938150198188
$ws.Range("B2:B25").Value = $colB

$colC = New-Object 'object[,]' 24,1
$colC[0,0] = 6.409948786382023
$colC[1,0] = 6.276082660054503
$colC[2,0] = 6.192341479013633
$colC[3,0] = 6.157864118017514
$colC[4,0] = 6.152119016962065
$colC[5,0] = 6.191877881272322
$colC[6,0] = 6.364134366419948
$colC[7,0] = 6.688125152107801
$colC[8,0] = 6.915827892633962
$colC[9,0] = 7.016805053282163
$colC[10,0] = 7.054639808124223
$colC[11,0] = 7.046509774908443
$colC[12,0] = 7.019925975386538
$colC[13,0] = 7.003589327381624
$colC[14,0] = 6.909173887248891
$colC[15,0] = 6.850564788863569
$colC[16,0] = 6.816610666629735
$colC[17,0] = 6.805073423256512
$colC[18,0] = 6.856829263444594
$colC[19,0] = 7.027745437813643
$colC[20,0] = 7.13708771344151
$colC[21,0] = 7.078954590461811
$colC[22,0] = 6.853997898969327
$colC[23,0] = 6.602160402581723
$ws.Range("C2:C25").Value = $colC

$colE = New-Object 'object[,]' 24,1
$colE[0,0] = 24.57431050371444
$colE[1,0] = 24.20900703491147
$colE[2,0] = 23.98743062499479
$colE[3,0] = 23.8979310382758
$colE[4,0] = 23.88312089377955
$colE[5,0] = 23.98622024006922
$colE[6,0] = 24.44785545905084
$colE[7,0] = 25.36984466089142
$colE[8,0] = 26.05093800629954
$colE[9,0] = 26.36024520352134
$colE[10,0] = 26.47718749755689
$colE[11,0] = 26.45201185191429
$colE[12,0] = 26.36987045428858
$colE[13,0] = 26.31952909693475
$colE[14,0] = 26.03070399021938
$colE[15,0] = 25.85330513762347
$colE[16,0] = 25.75122692060667
$colE[17,0] = 25.71666099968262
$colE[18,0] = 25.8721948162151
$colE[19,0] = 26.39400327104261
$colE[20,0] = 26.73390486969086
$colE[21,0] = 26.5526318254121
$colE[22,0] = 25.86365506206991
$colE[23,0] = 25.1193064554557
$ws.Range("E2:E25").Value = $colE

$colF = New-Object 'object[,]' 24,1
$colF[0,0] = 39.43079075513764
$colF[1,0] = 39.28267870295038
$colF[2,0] = 39.2029768627571
$colF[3,0] = 39.17334671822253
$colF[4,0] = 39.16859931354586
$colF[5,0] = 39.20256569779178
$colF[6,0] = 39.37740378421707
$colF[7,0] = 39.80825919335658
$colF[8,0] = 40.17661384442843
$colF[9,0] = 40.35496854584605
$colF[10,0] = 40.42401354753278
$colF[11,0] = 40.40907724878276
$colF[12,0] = 40.3606190091315
$colF[13,0] = 40.33113165733894
$colF[14,0] = 40.16517122354315
$colF[15,0] = 40.0660917110155
$colF[16,0] = 40.01012216655687
$colF[17,0] = 39.99134807907334
$colF[18,0] = 40.0765338157509
$colF[19,0] = 40.37481187468826
$colF[20,0] = 40.57850849371128
$colF[21,0] = 40.4690064937321
$colF[22,0] = 40.07180984486422
$colF[23,0] = 39.68245581113076
$ws.Range("F2:F25").Value = $colF

$colG = New-Object 'object[,]' 24,1
$colG[0,0] = 24.37199077321862
$colG[1,0] = 24.57488669666958
$colG[2,0] = 24.70866935439257
$colG[3,0] = 24.76549060526198
$colG[4,0] = 24.77506456512841
$colG[5,0] = 24.70942635180202
$colG[6,0] = 24.44003259045772
$colG[7,0] = 23.98526543937364
$colG[8,0] = 23.6966451713985
$colG[9,0] = 23.57537664052201
$colG[10,0] = 23.53091051618344
$colG[11,0] = 23.5404221273122
$colG[12,0] = 23.57168913191909
$colG[13,0] = 23.59103105941244
$colG[14,0] = 23.70477328718512
$colG[15,0] = 23.77712818907671
$colG[16,0] = 23.81968768421299
$colG[17,0] = 23.83425911978388
$colG[18,0] = 23.76932820340832
$colG[19,0] = 23.56246563694063
$colG[20,0] = 23.4357602332933
$colG[21,0] = 23.50260367111954
$colG[22,0] = 23.77285158230014
$colG[23,0] = 24.10035279992199
$ws.Range("G2:G25").Value = $colG

$colH = New-Object 'object[,]' 24,1
$colH[0,0] = 13.37699555961789
$colH[1,0] = 13.44073430844226
$colH[2,0] = 13.48209134234088
$colH[3,0] = 13.49950399047754
$colH[4,0] = 13.50242915250472
$colH[5,0] = 13.48232390968563
$colH[6,0] = 13.39851214547548
$colH[7,0] = 13.25174560485603
$colH[8,0] = 13.15458781356901
$colH[9,0] = 13.11269508904002
$colH[10,0] = 13.09716214629272
$colH[11,0] = 13.10049273337421
$colH[12,0] = 13.11141055512809
$colH[13,0] = 13.11814111745516
$colH[14,0] = 13.15737193664819
$colH[15,0] = 13.18202872649081
$colH[16,0] = 13.19642763120215
$colH[17,0] = 13.2013401340104
$colH[18,0] = 13.17938151717844
$colH[19,0] = 13.1081947490027
$colH[20,0] = 13.06359873598263
$colH[21,0] = 13.08722415659708
$colH[22,0] = 13.18057762461344
$colH[23,0] = 13.28957224688464
$ws.Range("H2:H25").Value = $colH

$colJ = New-Object 'object[,]' 24,1
$colJ[0,0] = 7.758437319948225
$colJ[1,0] = 7.782569894448106
$colJ[2,0] = 7.798197727053309
$colJ[3,0] = 7.804770476240586
$colJ[4,0] = 7.80587422947271
$colJ[5,0] = 7.79828554163794
$colJ[6,0] = 7.766590357378308
$colJ[7,0] = 7.710842277271221
$colJ[8,0] = 7.673756450190424
$colJ[9,0] = 7.6577189905957
$colJ[10,0] = 7.651765297610169
$colJ[11,0] = 7.653042231115406
$colJ[12,0] = 7.657226787393652
$colJ[13,0] = 7.659805477935127
$colJ[14,0] = 7.674821260816493
$colJ[15,0] = 7.684246004562207
$colJ[16,0] = 7.689745306782945
$colJ[17,0] = 7.691620761775487
$colJ[18,0] = 7.683234609642701
$colJ[19,0] = 7.655994446345087
$colJ[20,0] = 7.638886854104175
$colJ[21,0] = 7.647954013158633
$colJ[22,0] = 7.683691609277813
$colJ[23,0] = 7.725241231338076
$ws.Range("J2:J25").Value = $colJ

$colK = New-Object 'object[,]' 24,1
$colK[0,0] = 8.593915754289915
$colK[1,0] = 8.265763939214965
$colK[2,0] = 8.056141908174155
$colK[3,0] = 7.968757311034095
$colK[4,0] = 7.954131326954518
$colK[5,0] = 8.054971234196536
$colK[6,0] = 8.482504615459568
$colK[7,0] = 9.253173950067584
$colK[8,0] = 9.774501415901572
$colK[9,0] = 10.00132409122407
$colK[10,0] = 10.0856896425282
$colK[11,0] = 10.06758853457269
$colK[12,0] = 10.00829572110887
$colK[13,0] = 9.971777135852953
$colK[14,0] = 9.75946642230296
$colK[15,0] = 9.626542872321309
$colK[16,0] = 9.549118346716973
$colK[17,0] = 9.522738432330856
$colK[18,0] = 9.640793566255645
$colK[19,0] = 10.02575320843103
$colK[20,0] = 10.26842936514325
$colK[21,0] = 10.1397366431958
$colK[22,0] = 9.634353954716262
$colK[23,0] = 9.052342090700288
$ws.Range("K2:K25").Value = $colK

$colO = New-Object 'object[,]' 24,1
$colO[0,0] = 19.74384981818503
$colO[1,0] = 19.86558674690869
$colO[2,0] = 19.94490682421316
$colO[3,0] = 19.9783801314793
$colO[4,0] = 19.98400778482396
$colO[5,0] = 19.94535360196615
$colO[6,0] = 19.78487542133741
$colO[7,0] = 19.50647381342153
$colO[8,0] = 19.32407367389288
$colO[9,0] = 19.24590793128564
$colO[10,0] = 19.21700088879649
$colO[11,0] = 19.22319571418675
$colO[12,0] = 19.24351584354377
$colO[13,0] = 19.2560527444925
$colO[14,0] = 19.32927884999489
$colO[15,0] = 19.37543314355625
$colO[16,0] = 19.40243244571859
$colO[17,0] = 19.41165164358907
$colO[18,0] = 19.37047309291708
$colO[19,0] = 19.23752852198696
$colO[20,0] = 19.1546793084396
$colO[21,0] = 19.19852764847017
$colO[22,0] = 19.37271408471087
$colO[23,0] = 19.57790196034495
$ws.Range("O2:O25").Value = $colO
